# CU-1yrwt71 added UI tests for registration and code refactoring (#263)
#
# Adds "InvalidFirstName"/"InvalidLastName" columns (G/H) to the
# RegistrationTest sheet, fills in the remaining invalid-email/invalid-shelf
# test rows (including two brand-new rows), widens the affected columns,
# drops the now-stale hyperlinks on that sheet, and leaves the
# RegistrationTest tab as the active one (with F9 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationTest")

# --- Header row: two new "invalid" columns ---------------------------------
$ws.Range("G1").Value = "InvalidFirstName"
$ws.Range("H1").Value = "InvalidLastName"
$ws.Range("G1:H1").NumberFormat = "@"

# --- Row 2: long first/last name values for the two new columns ------------
$ws.Range("G2").Value = "StefanStefanStefanStefanStefanStefanStefanStefanSte"
$ws.Range("H2").Value = "GajicGajicGajicGajicGajicGajicGajicGajicGajicGajicG"
$ws.Range("G2:H2").NumberFormat = "@"

# --- Rows 3-5: blank placeholder cells for the two new columns -------------
$ws.Range("G3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("G4:H5").NumberFormat = "@"

# --- Row 6: the email column picks up a new (longer) invalid-email value ---
$ws.Range("E6").Value = "stefan.zgajic.stefan.zgajic.stefan.zgajic@gmail.com"
$ws.Range("G6:H6").NumberFormat = "@"

# --- Row 7 (new): one more invalid-shelf value ------------------------------
$ws.Range("F7").Value = "Shelf!@7Shelf!@7Shelf!@7Shelf!@7Shelf!@7Shelf!@7!@7"
$ws.Range("A7:E7").NumberFormat = "@"
$ws.Range("F7:H7").NumberFormat = "@"

# --- Row 8 (new): entirely blank placeholder row ----------------------------
$ws.Range("A8:H8").NumberFormat = "@"

# --- Drop the hyperlinks that used to decorate the sample data -------------
$ws.Hyperlinks.Delete()

# --- Column widths for the widened/new columns ------------------------------
$ws.Columns.Item(5).ColumnWidth = 43.17
$ws.Columns.Item(6).ColumnWidth = 51
$ws.Columns.Item(7).ColumnWidth = 45.67
$ws.Columns.Item(8).ColumnWidth = 42.33

# --- Make RegistrationTest the active tab with F9 selected -----------------
$ws.Activate()
$ws.Range("F9").Select()
